# New biosteam equipment lifetime feature:
# add a "Lifetime (yr)" row under the existing "Number" row (row 12) on
# Sheet1, mirroring its layout - label in column A, FALSE placeholders
# (same as the other per-equipment boolean rows) across columns B:AN.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Label cell, styled like the other left-hand category cells: left
# aligned with a thin left/right border (no top/bottom).
$label = $ws.Range("A13")
$label.Value = "Lifetime (yr)"
$label.HorizontalAlignment = -4131   # xlLeft
$label.Borders.Item(7).LineStyle = 1   # xlEdgeLeft  -> thin
$label.Borders.Item(10).LineStyle = 1  # xlEdgeRight -> thin

# Data cells: same FALSE/center-aligned look as row 12 (the "Number" row).
$data = $ws.Range("B13:AN13")
$data.Value = $false
$data.HorizontalAlignment = -4108    # xlCenter

# Move the view: scroll the unfrozen pane over a bit and select the new
# row's first cell (matches the author's view state after adding the row).
$win = $ws.Application.ActiveWindow
[void]($win.ScrollColumn = 28)
[void]($ws.Range("A13").Select())
